$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("SpMatrix")

$pairs = @(
    @(12, 17),
    @(13, 18),
    @(14, 19),
    @(15, 20),
    @(16, 21)
)

foreach ($pair in $pairs) {
    $r1 = $pair[0]
    $r2 = $pair[1]
    $rng1 = $ws.Range("A" + $r1 + ":FH" + $r1)
    $rng2 = $ws.Range("A" + $r2 + ":FH" + $r2)
    $v1 = $rng1.Value()
    $v2 = $rng2.Value()
    $rng1.Value = $v2
    $rng2.Value = $v1
}

$ws.Activate()
$ws.Range("K17").Select()
